$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("G2").Value = 0.1133093333333333
$ws.Range("H2").Value = 0.339928
$ws.Range("I2").Value = 0.02456654176752224
$ws.Range("J2").Value = 0.02456654176752224
$ws.Range("K2").Value = 3
$ws.Range("L2").Value = 1
$ws.Range("M2").Value = 2.731629
$ws.Range("N2").Value = 8.194887
$ws.Range("O2").Value = 0.5547800938501829
$ws.Range("P2").Value = 0.554780093850183
$ws.Range("Q2").Value = 0.309519060904
$ws.Range("R2").Value = 2.785671548136
$ws.Range("S2").Value = 0.01362902834736042
$ws.Range("T2").Value = 0.01362902834736043
$ws.Range("G3").Value = 0.1133093333333333
$ws.Range("H3").Value = 0.339928
$ws.Range("I3").Value = 0.02456654176752224
$ws.Range("J3").Value = 0.02456654176752224
$ws.Range("M3").Value = 0.06813733333333333
$ws.Range("O3").Value = 0.01383834927121065
$ws.Range("P3").Value = 0.01383834927121065
$ws.Range("Q3").Value = 0.00772059581511111
$ws.Range("R3").Value = 0.069485362336
$ws.Range("S3").Value = 0.0003399603853647572
$ws.Range("T3").Value = 0.0003399603853647573
$ws.Range("G4").Value = 0.1133093333333333
$ws.Range("H4").Value = 0.339928
$ws.Range("I4").Value = 0.02456654176752224
$ws.Range("J4").Value = 0.02456654176752224
$ws.Range("M4").Value = 2.124038666666666
$ws.Range("N4").Value = 6.372115999999999
$ws.Range("O4").Value = 0.4313815568786064
$ws.Range("P4").Value = 0.4313815568786064
$ws.Range("Q4").Value = 0.2406734052942222
$ws.Range("R4").Value = 2.166060647648
$ws.Range("S4").Value = 0.01059755303479705
$ws.Range("T4").Value = 0.01059755303479705
$ws.Range("I5").Value = 0.8380577451911468
$ws.Range("J5").Value = 0.8380577451911468
$ws.Range("K5").Value = 3
$ws.Range("L5").Value = 1
$ws.Range("M5").Value = 2.731629
$ws.Range("N5").Value = 8.194887
$ws.Range("O5").Value = 0.5547800938501829
$ws.Range("P5").Value = 0.554780093850183
$ws.Range("Q5").Value = 10.558866963433
$ws.Range("R5").Value = 95.02980267089698
$ws.Range("S5").Value = 0.4649377545290171
$ws.Range("T5").Value = 0.4649377545290172
$ws.Range("I6").Value = 0.8380577451911468
$ws.Range("J6").Value = 0.8380577451911468
$ws.Range("M6").Value = 0.06813733333333333
$ws.Range("O6").Value = 0.01383834927121065
$ws.Range("P6").Value = 0.01383834927121065
$ws.Range("Q6").Value = 0.2633787523524444
$ws.Range("S6").Value = 0.01159733578739834
$ws.Range("T6").Value = 0.01159733578739834
$ws.Range("I7").Value = 0.8380577451911468
$ws.Range("J7").Value = 0.8380577451911468
$ws.Range("M7").Value = 2.124038666666666
$ws.Range("N7").Value = 6.372115999999999
$ws.Range("O7").Value = 0.4313815568786064
$ws.Range("P7").Value = 0.4313815568786064
$ws.Range("Q7").Value = 8.210281010532887
$ws.Range("R7").Value = 73.89252909479599
$ws.Range("S7").Value = 0.3615226548747313
$ws.Range("T7").Value = 0.3615226548747313
$ws.Range("G8").Value = 0.6336240000000001
$ws.Range("H8").Value = 1.900872
$ws.Range("I8").Value = 0.1373757130413309
$ws.Range("J8").Value = 0.1373757130413309
$ws.Range("K8").Value = 3
$ws.Range("L8").Value = 1
$ws.Range("M8").Value = 2.731629
$ws.Range("N8").Value = 8.194887
$ws.Range("O8").Value = 0.5547800938501829
$ws.Range("P8").Value = 0.554780093850183
$ws.Range("Q8").Value = 1.730825693496
$ws.Range("R8").Value = 15.577431241464
$ws.Range("S8").Value = 0.07621331097380536
$ws.Range("T8").Value = 0.07621331097380538
$ws.Range("G9").Value = 0.6336240000000001
$ws.Range("H9").Value = 1.900872
$ws.Range("I9").Value = 0.1373757130413309
$ws.Range("J9").Value = 0.1373757130413309
$ws.Range("M9").Value = 0.06813733333333333
$ws.Range("O9").Value = 0.01383834927121065
$ws.Range("P9").Value = 0.01383834927121065
$ws.Range("Q9").Value = 0.043173449696
$ws.Range("R9").Value = 0.388561047264
$ws.Range("S9").Value = 0.001901053098447545
$ws.Range("T9").Value = 0.001901053098447545
$ws.Range("G10").Value = 0.6336240000000001
$ws.Range("H10").Value = 1.900872
$ws.Range("I10").Value = 0.1373757130413309
$ws.Range("J10").Value = 0.1373757130413309
$ws.Range("M10").Value = 2.124038666666666
$ws.Range("N10").Value = 6.372115999999999
$ws.Range("O10").Value = 0.4313815568786064
$ws.Range("P10").Value = 0.4313815568786064
$ws.Range("Q10").Value = 1.345841876128
$ws.Range("R10").Value = 12.112576885152
$ws.Range("S10").Value = 0.059261348969078
$ws.Range("T10").Value = 0.059261348969078
